# Weekly data update: insert two new rows of data at the top of the grape (Uva)
# price records block (rows 336-337), shifting the existing rows 336-419 down to 338-421.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at row 336 (existing rows 336:419 shift down to 338:421)
$ws.Range("A336:A337").EntireRow.Insert()

# --- New row 336: Flame Seedless, Provincia de Huasco ---
$ws.Cells.Item(336, 1).Value = 9
$ws.Cells.Item(336, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(336, 3).Value = "Metropolitana"
$ws.Cells.Item(336, 4).Value = 44551
$ws.Cells.Item(336, 5).Value = 13
$ws.Cells.Item(336, 6).Value = "Fruta"
$ws.Cells.Item(336, 7).Value = 100109
$ws.Cells.Item(336, 8).Value = "Uva"
$ws.Cells.Item(336, 9).Value = 100109001
$ws.Cells.Item(336, 10).Value = "Uva"
$ws.Cells.Item(336, 11).Value = "Flame Seedless"
$ws.Cells.Item(336, 12).Value = "Primera"
$ws.Cells.Item(336, 13).Value = 380
$ws.Cells.Item(336, 14).Value = 8000
$ws.Cells.Item(336, 15).Value = 8000
$ws.Cells.Item(336, 16).Value = 8000
$ws.Cells.Item(336, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(336, 18).Value = "Provincia de Huasco"
$ws.Cells.Item(336, 19).Value = 800
$ws.Cells.Item(336, 20).Value = 10

# --- New row 337: Superior Seedless, Provincia de Huasco ---
$ws.Cells.Item(337, 1).Value = 9
$ws.Cells.Item(337, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(337, 3).Value = "Metropolitana"
$ws.Cells.Item(337, 4).Value = 44551
$ws.Cells.Item(337, 5).Value = 13
$ws.Cells.Item(337, 6).Value = "Fruta"
$ws.Cells.Item(337, 7).Value = 100109
$ws.Cells.Item(337, 8).Value = "Uva"
$ws.Cells.Item(337, 9).Value = 100109001
$ws.Cells.Item(337, 10).Value = "Uva"
$ws.Cells.Item(337, 11).Value = "Superior Seedless"
$ws.Cells.Item(337, 12).Value = "Primera"
$ws.Cells.Item(337, 13).Value = 350
$ws.Cells.Item(337, 14).Value = 16000
$ws.Cells.Item(337, 15).Value = 16000
$ws.Cells.Item(337, 16).Value = 16000
$ws.Cells.Item(337, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(337, 18).Value = "Provincia de Huasco"
$ws.Cells.Item(337, 19).Value = 1600
$ws.Cells.Item(337, 20).Value = 10
